$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("J1:J2").EntireColumn.Insert()
$ws.Range("J1").Value = "In Zahlungslauf" + [char]10 + "ignorieren"
$ws.Range("J2").Value = "{isIgnoriert}"
$ws.Columns.Item(9).ColumnWidth = 12.592447916666666
$ws.Columns.Item(10).ColumnWidth = 15.307291666666666
Write-Output "done"
